$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Containers" worksheet right before the "Enum" sheet.
# ---------------------------------------------------------------------------
$enumSheet = $wb.Worksheets.Item("Enum")
$propsSheet = $wb.Worksheets.Item("Properties")
$viewsSheet = $wb.Worksheets.Item("Views")

$containers = $wb.Worksheets.Add($enumSheet)
$containers.Name = "Containers"

# Reuse the existing "title row" formatting (bold size-20 font on the orange
# fill, fontId 2 / fillId 2) from another sheet's A1 so no new font/fill is
# created - then tweak the alignment to match the target (left for the text
# cell, general for the filler cells).
$viewsSheet.Range("A1").Copy()
$containers.Range("A1:F1").PasteSpecial(-4122)
$containers.Range("A1").HorizontalAlignment = -4131
$containers.Range("B1:F1").HorizontalAlignment = 1

# Header row re-uses the existing bold sub-header font (fontId 3).
$propsSheet.Range("A2").Copy()
$containers.Range("A2:F2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# New text values are entered in the same order the original author typed
# them (this controls the order new entries are appended to the shared
# string table), before the formatting/header reuse above is touched again.
# ---------------------------------------------------------------------------
$propsSheet.Range("B60").Value = "classicEquipmentGUID"
$containers.Range("A1").Value = "Definition of Containers"
$containers.Range("E2").Value = "Used For"
$containers.Range("F2").Value = "Neat ID"
$containers.Range("F3").Value = "http://purl.org/cognite/neat/neatId_bb7e4121_bc17_46c7_a93d_c855f710ccbe"
$containers.Range("E3").Value = "node"

# Remaining header cells / data reuse already-existing shared strings.
$containers.Range("A2").Value = "Container"
$containers.Range("B2").Value = "Name"
$containers.Range("C2").Value = "Description"
$containers.Range("D2").Value = "Constraint"
$containers.Range("A3").Value = "ClassicEquipment"

$containers.Rows.Item(1).RowHeight = 25.8
$containers.Rows.Item(2).RowHeight = 18

$containers.Columns.Item(1).ColumnWidth = 38.5546875
$containers.Columns.Item(5).ColumnWidth = 10.33203125
$containers.Columns.Item(6).ColumnWidth = 67.6640625

$containers.Activate()
$containers.Range("E8").Select()

# ---------------------------------------------------------------------------
# 2. Properties sheet: add the classicEquipmentGUID mapping row (row 60) and
#    update the selection / scroll position.
# ---------------------------------------------------------------------------
$propsSheet.Range("A60").Value = "ClassicEquipment"
$propsSheet.Range("F60").Value = "text"
$propsSheet.Range("G60").Value = $true
$propsSheet.Range("H60").Value = $false
$propsSheet.Range("I60").Value = $false
$propsSheet.Range("K60").Value = "ClassicEquipment"
$propsSheet.Range("L60").Value = "classicEquipmentGUID"

$propsSheet.Activate()
$propsSheet.Range("A39").Select()
$propsSheet.Range("K66").Select()

# ---------------------------------------------------------------------------
# 3. Views sheet: update the stored selection to the A1:G10 range.
# ---------------------------------------------------------------------------
$viewsSheet.Activate()
$viewsSheet.Range("A1:G10").Select()

# ---------------------------------------------------------------------------
# 4. Leave "Properties" as the active sheet/tab (matches tabSelected moving
#    from "Metadata" to "Properties").
# ---------------------------------------------------------------------------
$propsSheet.Activate()
$propsSheet.Range("K66").Select()
